$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.81 = 6555.39 pesos`n✅ 6555.39 pesos = 1.8 = 948.26 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the updated rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("O10").Value = 3626.05
$wsTasas.Range("N12").Value = 3643.96
$wsTasas.Range("O12").Value = 527.111
